$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.633.14"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "3.094.85"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "516.04"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.19"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.08%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.30"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  -0.47%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.374"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "3.623.22"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("E13").Value = "  +2.40%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "25.72"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.30%  "
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "57.723.17"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "3.088.06"
$ws.Range("E17").Value = "  +0.65%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.13"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.66%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.14"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("E20").Value = "  +0.02%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "336.68"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("E22").Value = "  +0.12%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.504"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.84"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("E27").Value = "  +2.95%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.46"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.41%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.10"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  -3.62%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "153.41"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.53"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "27.28"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +6.63%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.91"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  -3.17%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0684"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "3.136.96"
$ws.Range("E39").Value = "  +1.15%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "36.87"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.672"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.86"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "2.292.94"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("E45").Value = "  -0.61%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0252"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "20.31"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.947"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("E50").Value = "  +1.09%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.689"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.67%  "
